# Auto-generated Excel COM-interop edit script
# Applies cell value updates to the active sheet (cryptos list) as described in the diff.
# Values that look like plain numbers must be forced to Text format (NumberFormat "@")
# before assignment, otherwise Excel auto-converts them to numeric cells (losing the
# original text representation / introducing floating-point noise), whereas the source
# workbook stores every Price/Volume cell as literal text (inlineStr).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.817.89'
$ws.Range('E2').Value = '  +0.37%  '
$ws.Range('D3').Value = '2.306.18'
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '301.62'
$ws.Range('E5').Value = '  -1.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.04'
$ws.Range('E6').Value = '  -0.52%  '
$ws.Range('E7').Value = '  +0.46%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.494'
$ws.Range('E9').Value = '  -1.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.65'
$ws.Range('E10').Value = '  -2.59%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.19'
$ws.Range('E11').Value = '  +4.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0789'
$ws.Range('E12').Value = '  -0.14%  '
$ws.Range('E13').Value = '  -0.20%  '
$ws.Range('E14').Value = '  +0.72%  '
$ws.Range('D15').Value = '2.669.22'
$ws.Range('E15').Value = '  +1.33%  '
$ws.Range('D16').Value = '2.304.07'
$ws.Range('E16').Value = '  +0.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.785'
$ws.Range('E17').Value = '  +0.76%  '
$ws.Range('D18').Value = '42.753.21'
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.21'
$ws.Range('E19').Value = '  -5.48%  '
$ws.Range('D20').Value = '0.0₃0891'
$ws.Range('E20').Value = '  -0.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.02'
$ws.Range('E21').Value = '  +0.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.79'
$ws.Range('E22').Value = '  +1.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.28'
$ws.Range('E23').Value = '  +7.53%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '235.15'
$ws.Range('E24').Value = '  -0.42%  '
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.42'
$ws.Range('E26').Value = '  -1.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.43'
$ws.Range('E27').Value = '  -2.89%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.36'
$ws.Range('E28').Value = '  +14.72%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '165.68'
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.07'
$ws.Range('E30').Value = '  +0.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.07'
$ws.Range('E31').Value = '  -3.20%  '
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('E33').Value = '  +0.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.62'
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.45'
$ws.Range('E35').Value = '  -6.79%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0702'
$ws.Range('E36').Value = '  +1.67%  '
$ws.Range('E37').Value = '  -2.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.100'
$ws.Range('E38').Value = '  -0.81%  '
$ws.Range('E39').Value = '  +0.57%  '
$ws.Range('E40').Value = '  -0.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.70'
$ws.Range('E41').Value = '  +0.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '20.09'
$ws.Range('E42').Value = '  +10.88%  '
$ws.Range('D43').Value = '1.970.67'
$ws.Range('E43').Value = '  -1.53%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.48'
$ws.Range('E44').Value = '  +5.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0279'
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('E46').Value = '  -2.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.77'
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('D48').Value = '2.533.58'
$ws.Range('E48').Value = '  +1.24%  '
$ws.Range('E49').Value = '  -0.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '53.36'
$ws.Range('E50').Value = '  -0.31%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '71.49'
$ws.Range('E51').Value = '  +0.14%  '
